$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.458.71"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "3.907.26"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +9.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.721"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("E10").Value = "  -5.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000335"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.18"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "4.533.48"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").Value = "3.896.05"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.08"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.23"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +8.46%  "
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.76"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("D20").Value = "69.383.41"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "429.75"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.22"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.61"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.02"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.59"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.41"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "686.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.19"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.34%  "
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "67.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +11.57%  "
$ws.Range("E34").Value = "  +12.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.99"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").Value = "0.0₃0848"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.997"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.18"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.94%  "
$ws.Range("E43").Value = "  +6.01%  "
$ws.Range("E44").Value = "  -5.06%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.140"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.34"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "0.0₆0357"
$ws.Range("E47").Value = "  +12.49%  "
$ws.Range("E48").Value = "  +7.58%  "
$ws.Range("D49").Value = "2.746.99"
$ws.Range("E49").Value = "  +13.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.43"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.27"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.04%  "
